$d = $word.ActiveDocument

function ReplaceInPara($idx, $find, $repl) {
    $p = $d.Paragraphs.Item($idx)
    $r = $d.Range($p.Range.Start, $p.Range.End)
    $found = $r.Find.Execute($find, $true, $false, $false, $false, $false, $true, 0, $false, $repl, 1)
    if (-not $found) {
        Write-Output "NOT FOUND: idx=$idx find=$find"
    }
}

ReplaceInPara 1 'How to run your coding club safely' 'Πώς να λειτουργήσετε με ασφάλεια τη λέσχη προγραμματισμού σας'
ReplaceInPara 3 'This table helps you identify your options for running your club while your community is still at risk from coronavirus' 'Αυτός ο πίνακας σάς βοηθά να εντοπίσετε τις επιλογές σας για τη λειτουργία της λέσχης σας ενώ η κοινότητά σας εξακολουθεί να κινδυνεύει από τον κορονοϊό'
ReplaceInPara 5 'You should:' 'Θα πρέπει:'
ReplaceInPara 7 'Always follow your local public health advice, consider the risks in your location, and decide whether to run your club sessions in  person in consultation with the venue that hosts it.' 'Ακολουθείτε πάντα τις τοπικές συμβουλές για τη δημόσια υγεία, εξετάζετε τους κινδύνους στην τοποθεσία σας και αποφασίζετε εάν θα πραγματοποιήσετε τις συνεδρίες της λέσχης σας αυτοπροσώπως σε συνεννόηση με τον χώρο που τη φιλοξενεί.'
ReplaceInPara 8 'Regularly check for updates to your local public health advice.' 'Ελέγχετε τακτικά για ενημερώσεις στις τοπικές συμβουλές για τη δημόσια υγεία.'
ReplaceInPara 9 'Follow our ' 'Ακολουθείτε την '
ReplaceInPara 9 'safeguarding policy' 'πολιτική προστασίας μας'
ReplaceInPara 9 ', and make sure you attend one of our ‘Best practices for online sessions’ community calls.' ' και βεβαιωθείτε ότι ότι θα παρακολουθήσετε μια από τις προσκλήσεις της κοινότητας «Βέλτιστες πρακτικές για διαδικτυακές συνεδρίες».'
ReplaceInPara 12 'My current situation' 'Η τρέχουσα κατάσταση μου'
ReplaceInPara 13 'Recommended activities' 'Προτεινόμενες δραστηριότητες'
ReplaceInPara 14 'What should I do next?' 'Τι πρέπει να κάνω μετά;'
ReplaceInPara 16 'It is ' ''
ReplaceInPara 16 'not safe to run in-person club sessions' 'Δεν είναι ασφαλές να πραγματοποιείτε δια ζώσης συνεδρίες'
ReplaceInPara 16 ' right now.' ' αυτήν τη στιγμή.'
ReplaceInPara 17 'Online sessions' 'Online συνεδρίες'
ReplaceInPara 18 'Remote activities' 'Απομακρυσμένες δραστηριότητες'
ReplaceInPara 19 'A mixture of both' 'Ένας συνδυασμός και των δύο'
ReplaceInPara 20 'Read our guidance on how to run ' 'Διαβάστε τις οδηγίες μας σχετικά με τον τρόπο λειτουργίας '
ReplaceInPara 20 'online sessions' 'online συνεδριών'
ReplaceInPara 20 ' and ' ' και '
ReplaceInPara 20 'remote activities' 'απομακρυσμένων δραστηριοτήτων'
ReplaceInPara 24 'I am able to run ' 'Είμαι σε θέση να οργανώσω '
ReplaceInPara 24 'in-person club sessions with safety measures' 'δια ζώσης συνεδρίες της λέσχης με μέτρα ασφαλείας'
ReplaceInPara 24 ' such as social distancing.' ', όπως η κοινωνική απόσταση.'
ReplaceInPara 25 'In-person sessions respecting  local safety guidance' 'Δια ζώσης συνεδρίες με τήρηση των τοπικών κανόνων ασφαλείας'
ReplaceInPara 26 'Online sessions' 'Online συνεδρίες'
ReplaceInPara 27 'Remote activities' 'Απομακρυσμένες δραστηριότητες'
ReplaceInPara 28 'A mixture of two or more of  the above' 'Ένας συνδυασμός δύο ή περισσότερων από τα παραπάνω'
ReplaceInPara 29 'Read our guidance on how to run ' 'Διαβάστε τις οδηγίες μας σχετικά με τον τρόπο λειτουργίας '
ReplaceInPara 29 'online sessions' 'online συνεδριών'
ReplaceInPara 29 ' and ' ' και '
ReplaceInPara 29 'remote activities' 'απομακρυσμένων δραστηριοτήτων'
ReplaceInPara 31 'For in-person sessions: ' 'Για δια ζώσης συνεδρίες: '
ReplaceInPara 32 'Plan your sessions in line with the safety guidelines of your local authorities  and host venue to make sure you are prepared.' 'Προγραμματίστε τις συνεδρίες σας σύμφωνα με τις οδηγίες ασφαλείας των τοπικών αρχών και του χώρου φιλοξενίας για να βεβαιωθείτε ότι είστε προετοιμασμένοι.'
ReplaceInPara 33 'Ensure equal opportunities for young people and adults at increased risk from  coronavirus and offer them safe activities, or alternatives to meeting in person.' 'Εξασφαλίστε ίσες ευκαιρίες για νέους και ενήλικες που διατρέχουν αυξημένο κίνδυνο από τον κορωνοϊό και προσφέρετέ τους ασφαλείς δραστηριότητες ή εναλλακτικές λύσεις αντί για δια ζώσης συνάντηση.'
ReplaceInPara 34 'In case safety guidelines in your local area change, prepare by training all your  participants to use the tools needed to access your online sessions or  participate in remote activities.' 'Σε περίπτωση που αλλάξουν οι οδηγίες ασφαλείας στην περιοχή σας, προετοιμαστείτε εκπαιδεύοντας όλους τους συμμετέχοντες να χρησιμοποιούν τα εργαλεία που χρειάζονται για πρόσβαση στις online συνεδρίες σας ή για συμμετοχή σε απομακρυσμένες δραστηριότητες.'
ReplaceInPara 36 'I am able to ' 'Είμαι σε θέση να '
ReplaceInPara 36 'run in-person club sessions with no safety measures.' 'πραγματοποιώ δια ζώσης συνεδρίες της λέσχης χωρίς μέτρα ασφαλείας'
ReplaceInPara 37 'In-person sessions' 'Δια ζώσης συνεδρίες'
ReplaceInPara 38 'A mixture of in-person sessions,  online sessions, and remote  activities' 'Ένας συνδυασμός από δια ζώσης συνεδρίες, online συνεδρίες και απομακρυσμένες δραστηριότητες'
ReplaceInPara 39 'Ensure equal opportunities for young people and adults at increased risk from coronavirus and offer them safe activities, or alternatives to meeting in person.' 'Εξασφαλίστε ίσες ευκαιρίες για νέους και ενήλικες που διατρέχουν αυξημένο κίνδυνο από τον κορωνοϊό και προσφέρετέ τους ασφαλείς δραστηριότητες ή εναλλακτικές λύσεις αντί για δια ζώσης συνάντηση.'
ReplaceInPara 40 'Read our guidance on how to run ' 'Διαβάστε τις οδηγίες μας σχετικά με τον τρόπο λειτουργίας '
ReplaceInPara 40 'online sessions' 'online συνεδριών'
ReplaceInPara 40 ' and ' ' και '
ReplaceInPara 40 'remote activities' 'απομακρυσμένων δραστηριοτήτων'
ReplaceInPara 67 'Code Club and CoderDojo are part of the Raspberry Pi Foundation, UK registered charity 1129409 ' 'Το Code Club και το CoderDojo αποτελούν μέρος του Raspberry Pi Foundation, εγγεγραμμένης φιλανθρωπικής οργάνωσης με έδρα το Ηνωμένο Βασίλειο (1129409) '
